# Jogos_da_Semana_FlashScore_2024-11-07.xlsx - apply weekly odds refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds in rows 2-5 (refreshed FlashScore prices) ---
$ws.Range("V2").Value = 1.58
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.77
$ws.Range("U4").Value = 2.62
$ws.Range("V4").Value = 1.41
$ws.Range("G5").Value = 1.45
$ws.Range("I5").Value = 7.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 2.05
$ws.Range("V5").Value = 1.7
$ws.Range("W5").Value = 6.5
$ws.Range("X5").Value = 6.5
$ws.Range("Y5").Value = 8.5
$ws.Range("Z5").Value = 9.5
$ws.Range("AB5").Value = 29
$ws.Range("AD5").Value = 8
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 451
$ws.Range("AH5").Value = 17
$ws.Range("AR5").Value = 41
$ws.Range("AT5").Value = 3
$ws.Range("AW5").Value = 8
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 351
$ws.Range("BC5").Value = 151

# --- Append new row 8: MEXICO - LIGA DE EXPANSION MX, Tepatitlan de Morelos vs Tapatio ---
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "Q1OUVsW1"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "07/11/2024"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "22:00"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "Tepatitlan de Morelos"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Tapatio"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = 2.67
$ws.Range("H8").Value = 2.92
$ws.Range("I8").Value = 2.67
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 2.02
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 7.6
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 2.82
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.5
$ws.Range("U8").Value = 1.7
$ws.Range("V8").Value = 1.93
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 13.5
$ws.Range("Y8").Value = 9.75
$ws.Range("Z8").Value = 32
$ws.Range("AA8").Value = 23
$ws.Range("AB8").Value = 32
$ws.Range("AC8").Value = 8.5
$ws.Range("AD8").Value = 5.7
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 450
$ws.Range("AH8").Value = 8.5
$ws.Range("AI8").Value = 14
$ws.Range("AJ8").Value = 9.75
$ws.Range("AK8").Value = 32
$ws.Range("AL8").Value = 23
$ws.Range("AM8").Value = 30
$ws.Range("AN8").Value = 4.6
$ws.Range("AO8").Value = 14
$ws.Range("AP8").Value = 20
$ws.Range("AQ8").Value = 60
$ws.Range("AR8").Value = 90
$ws.Range("AS8").Value = 250
$ws.Range("AT8").Value = 2.55
$ws.Range("AU8").Value = 6.3
$ws.Range("AV8").Value = 50
$ws.Range("AW8").Value = 4.65
$ws.Range("AX8").Value = 14
$ws.Range("AY8").Value = 19
$ws.Range("AZ8").Value = 60
$ws.Range("BA8").Value = 80
$ws.Range("BB8").Value = 200
$ws.Range("BC8").Value = 51
$ws.Range("BD8").Value = 51
